$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
$sm = $d.SlideMaster
Write-Host "SlideMaster Name:" $sm.Name
$th = $sm.Theme
Write-Host "Theme:" $th
Write-Host "Theme.Name:" $th.Name
